$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3098
$ws1.Range("F3").Value = 505
$ws1.Range("F6").Value = 16
$ws1.Range("F8").Value = 20
$ws1.Range("F9").Value = 1080
$ws1.Range("F10").Value = 15153
$ws1.Range("F11").Value = 201
$ws1.Range("F12").Value = 151
$ws1.Range("F13").Value = 512
$ws1.Range("F14").Value = 6012
$ws1.Range("F15").Value = 612
$ws1.Range("F17").Value = 56
$ws1.Range("F21").Value = 105
$ws1.Range("F24").Value = 836
$ws1.Range("F25").Value = 4965
$ws1.Range("F26").Value = 102
$ws1.Range("F27").Value = 10850
$ws1.Range("F29").Value = 3
$ws1.Range("F31").Value = 141
$ws1.Range("F32").Value = 3772

# Sheet: 演出 (Show)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 15

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3098
$ws4.Range("F4").Value = 505
$ws4.Range("F7").Value = 16
$ws4.Range("F9").Value = 20
$ws4.Range("F10").Value = 1080
$ws4.Range("F11").Value = 15153
$ws4.Range("F12").Value = 201
$ws4.Range("F13").Value = 151
$ws4.Range("F14").Value = 512
$ws4.Range("F15").Value = 6012
$ws4.Range("F16").Value = 612
$ws4.Range("F18").Value = 56
$ws4.Range("F22").Value = 105
$ws4.Range("F25").Value = 836
$ws4.Range("F26").Value = 4965
$ws4.Range("F27").Value = 102
$ws4.Range("F28").Value = 15
$ws4.Range("F29").Value = 10850
$ws4.Range("F31").Value = 3
$ws4.Range("F33").Value = 141
$ws4.Range("F34").Value = 3772
